$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 17, shifting existing rows (17-42) down to (18-43).
$ws.Rows("17:17").Insert()

# Populate the newly inserted row 17 with the new weekly record.
$ws.Cells.Item(17, 1).Value = 9
$ws.Cells.Item(17, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(17, 3).Value = "Metropolitana"
$ws.Cells.Item(17, 4).Value = 44973
$ws.Cells.Item(17, 5).Value = 13
$ws.Cells.Item(17, 6).Value = 100112010
$ws.Cells.Item(17, 7).Value = "Achicoria"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 90
$ws.Cells.Item(17, 11).Value = 7000
$ws.Cells.Item(17, 12).Value = 8000
$ws.Cells.Item(17, 13).Value = 7500
$ws.Cells.Item(17, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(17, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(17, 16).Value = 469
$ws.Cells.Item(17, 17).Value = 16
$ws.Cells.Item(17, 18).Value = "Hortaliza"

# The Origen values of the two rows that used to be 36 and 37 (now 37 and 38)
# swap which one reads "Región Metropolitana" vs "Provincia de Quillota" as a
# side effect of the insert/shift in the source data; make sure they land
# exactly as in the target sheet.
$ws.Cells.Item(37, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(38, 15).Value = "Región Metropolitana"
